$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Participant ID query text (replaces the old one previously used in B2)
$newParticipantQuery = "MATCH (p:participant)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nOPTIONAL MATCH (p)<--(diag:diagnosis)`nOPTIONAL MATCH (samp)<--(f:file)`nOPTIONAL MATCH (f)<--(g:genomic_info)`nWITH s, p, samp, f, g, diag`nWHERE g.instrument_model in ['Illumina NextSeq']`nWITH p`nOPTIONAL MATCH (p)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nWITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp`nRETURN `ncoalesce(p.participant_id,'') as ``Participant ID``,`ncoalesce(s.study_name, '') as ``Study Name``,`ncoalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(p.gender,'') as ``Gender``,`ncoalesce(apoc.text.join(samp, ','), '') as ``Samples```nORDER BY p.participant_id limit 100"

# Existing query texts that are reused but now live in a different cell
$sampleIdQuery = "Match (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina NextSeq']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`nRETURN  `n coalesce(samp.sample_id, '') as ``Sample ID``,`n coalesce(p.participant_id,'') as ``Participant ID``,`n coalesce(s.study_name, '') as ``Study Name``,`n coalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(samp.sample_tumor_status,'') as ``Tumor``,`ncoalesce(samp.sample_type,'') as ``Analyte Type```nORDER By samp.sample_id LIMIT 100"

$fileNameQuery = "Match (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina NextSeq']MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`nRETURN `n    coalesce(f.file_name, '') as ``File Name``,`n    coalesce(s.study_name, '') as ``Study Name``,`n    coalesce(s.phs_accession,'') as ``Accession``,`n    coalesce(p.participant_id,'') as ``Participant ID``,`n    coalesce(samp.sample_id, '') as ``Sample ID``,`n    coalesce(f.file_type, '') as ``File Type```nORDER By f.file_name LIMIT 100"

$statQuery = "MATCH (f:file)`nMatch (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina NextSeq']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,f, s, collect(distinct samp.sample_id) as samp`nRETURN`ncount(distinct s) AS Studies,`ncount(distinct p) AS Participants,`ncount(distinct samp) AS Samples,`ncount(distinct f) AS Files"

# Row 2 (SamplesTab)
$ws.Range("B2").Value = $newParticipantQuery
$ws.Range("C2").Value = $statQuery

# Row 3 (FilesTab)
$ws.Range("B3").Value = $sampleIdQuery
$ws.Range("C3").Value = $statQuery

# Row 4 (ParticipantsTab)
$ws.Range("B4").Value = $fileNameQuery
$ws.Range("C4").Value = $statQuery

# Update selection to C2
$ws.Range("C2").Select()
